# Ready for Part D insertion
# Adds a "flowDir" column (AD) set to "Countercurrent" for every data row,
# and corrects the Tci/Tco (columns H/I) values, which had been swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column AD: "flowDir" header + "Countercurrent" for each data row ---
# Copy the existing header formatting (A1) onto the new header cell so the
# new column matches the look of the rest of row 1.
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("AD1").Value = "flowDir"

$ws.Range("A1").Copy($ws.Range("AD2:AD6"))
$ws.Range("AD2:AD6").Value = "Countercurrent"

$ws.Range("AD1").EntireColumn.ColumnWidth = 13.9

# --- Fix Tci (H) / Tco (I) values: they were swapped for every run ---
$hCol = "H"
$iCol = "I"
for ($row = 2; $row -le 6; $row++) {
    $hCell = $ws.Range("$hCol$row")
    $iCell = $ws.Range("$iCol$row")
    $hVal = $hCell.Value()
    $iVal = $iCell.Value()
    $hCell.Value = $iVal
    $iCell.Value = $hVal
}

# Update the used range / selection to reflect the new column, matching the
# post-edit workbook state.
[void]$ws.Range("K7").Select()
